$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab05")

# --- Header text corrections (row 2, columns J and K) ---
$ws.Range("J2").Value = "Taux de dépendance des personnes âgées 65+/(15-64), 2022*100, 2022"
$ws.Range("K2").Value = "Taux de dépendance des enfants à charge <15/(15-64), 2022*100, 2022"

# --- Updated aggregate figures for "Afrique, États fragiles" (row 97) ---
$ws.Range("C97").Value = 842549.25199999998
$ws.Range("D97").Value = 348533.71888524003
$ws.Range("E97").Value = 494015.53311476001
$ws.Range("F97").Value = 48796.272602885401
$ws.Range("G97").Value = 85
$ws.Range("H97").Value = 117932.369893483
$ws.Range("I97").Value = 39
$ws.Range("J97").Value = 5.4141204157620599
$ws.Range("K97").Value = 72.908692094416395
$ws.Range("L97").Value = 78.322812510178395

# --- Updated aggregate figures for "RDM, États fragiles" (row 98) ---
$ws.Range("C98").Value = 692226.44200000004
$ws.Range("D98").Value = 318611.69368165999
$ws.Range("E98").Value = 373614.74831833999
$ws.Range("F98").Value = 36149.9535672317
$ws.Range("G98").Value = 62
$ws.Range("H98").Value = 129607.049632435
$ws.Range("I98").Value = 33
$ws.Range("J98").Value = 8.7123928268905999
$ws.Range("K98").Value = 47.823885843296999
$ws.Range("L98").Value = 56.5362786701875

# --- Minor recalculated rounding of L72 ---
$ws.Range("L72").Value = 78.326488430312494
